$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 144.783305
$ws.Range("H2").Value = 434.349915
$ws.Range("I2").Value = 0.2430046335191003
$ws.Range("J2").Value = 0.251012682214973
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.2385586666666667
$ws.Range("N2").Value = 0.715676
$ws.Range("O2").Value = 0.003173413598159477
$ws.Range("P2").Value = 0.003191469360920207
$ws.Range("Q2").Value = 34.53931219639333
$ws.Range("R2").Value = 310.85380976754
$ws.Range("S2").Value = 0.0007711542084252733
$ws.Range("T2").Value = 0.0008010992844914868

$ws.Range("G3").Value = 144.783305
$ws.Range("H3").Value = 434.349915
$ws.Range("I3").Value = 0.2430046335191003
$ws.Range("J3").Value = 0.251012682214973
$ws.Range("M3").Value = 28.36841766666667
$ws.Range("N3").Value = 85.105253
$ws.Range("O3").Value = 0.377369322353974
$ws.Range("P3").Value = 0.3795164395660363
$ws.Range("Q3").Value = 4107.273267400388
$ws.Range("R3").Value = 36965.45940660349
$ws.Range("S3").Value = 0.09170249387997868
$ws.Range("T3").Value = 0.09526343944014745

$ws.Range("G4").Value = 144.783305
$ws.Range("H4").Value = 434.349915
$ws.Range("I4").Value = 0.2430046335191003
$ws.Range("J4").Value = 0.251012682214973
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.5
$ws.Range("M4").Value = 1.275895
$ws.Range("N4").Value = 2.55179
$ws.Range("O4").Value = 0.01697252336039082
$ws.Range("P4").Value = 0.01137939458707931
$ws.Range("Q4").Value = 184.728294932975
$ws.Range("R4").Value = 1108.36976959785
$ws.Range("S4").Value = 0.00412440181908614
$ws.Range("T4").Value = 0.002856372357285323

$ws.Range("G5").Value = 144.783305
$ws.Range("H5").Value = 434.349915
$ws.Range("I5").Value = 0.2430046335191003
$ws.Range("J5").Value = 0.251012682214973
$ws.Range("M5").Value = 45.29127766666667
$ws.Range("N5").Value = 135.873833
$ws.Range("O5").Value = 0.6024847406874758
$ws.Range("P5").Value = 0.6059126964859642
$ws.Range("Q5").Value = 6557.420868252688
$ws.Range("R5").Value = 59016.78781427418
$ws.Range("S5").Value = 0.1464065836116103
$ws.Range("T5").Value = 0.1520917711330487

$ws.Range("G6").Value = 82.248871
$ws.Range("I6").Value = 0.1380466950572427
$ws.Range("J6").Value = 0.1425959278859072
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.2385586666666667
$ws.Range("N6").Value = 0.715676
$ws.Range("O6").Value = 0.003173413598159477
$ws.Range("P6").Value = 0.003191469360920207
$ws.Range("Q6").Value = 19.62118100059866
$ws.Range("R6").Value = 176.590629005388
$ws.Range("S6").Value = 0.0004380792592756286
$ws.Range("T6").Value = 0.0004550905348398602

$ws.Range("G7").Value = 82.248871
$ws.Range("I7").Value = 0.1380466950572427
$ws.Range("J7").Value = 0.1425959278859072
$ws.Range("M7").Value = 28.36841766666667
$ws.Range("N7").Value = 85.105253
$ws.Range("O7").Value = 0.377369322353974
$ws.Range("P7").Value = 0.3795164395660363
$ws.Range("Q7").Value = 2333.270325139788
$ws.Range("R7").Value = 20999.43292625809
$ws.Range("S7").Value = 0.05209458776695736
$ws.Range("T7").Value = 0.05411749884787476

$ws.Range("G8").Value = 82.248871
$ws.Range("I8").Value = 0.1380466950572427
$ws.Range("J8").Value = 0.1425959278859072
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.5
$ws.Range("M8").Value = 1.275895
$ws.Range("N8").Value = 2.55179
$ws.Range("O8").Value = 0.01697252336039082
$ws.Range("P8").Value = 0.01137939458707931
$ws.Range("Q8").Value = 104.940923264545
$ws.Range("R8").Value = 629.64553958727
$ws.Range("S8").Value = 0.002343000756683799
$ws.Range("T8").Value = 0.001622655329924445

$ws.Range("G9").Value = 82.248871
$ws.Range("I9").Value = 0.1380466950572427
$ws.Range("J9").Value = 0.1425959278859072
$ws.Range("M9").Value = 45.29127766666667
$ws.Range("N9").Value = 135.873833
$ws.Range("O9").Value = 0.6024847406874758
$ws.Range("P9").Value = 0.6059126964859642
$ws.Range("Q9").Value = 3725.156454230847
$ws.Range("R9").Value = 33526.40808807762
$ws.Range("S9").Value = 0.0831710272743259
$ws.Range("T9").Value = 0.08640068317326813

$ws.Range("G10").Value = 163.8590903333333
$ws.Range("H10").Value = 491.577271
$ws.Range("I10").Value = 0.2750214756820535
$ws.Range("J10").Value = 0.284084617144743
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.2385586666666667
$ws.Range("N10").Value = 0.715676
$ws.Range("O10").Value = 0.003173413598159477
$ws.Range("P10").Value = 0.003191469360920207
$ws.Range("Q10").Value = 39.09000611113289
$ws.Range("R10").Value = 351.810055000196
$ws.Range("S10").Value = 0.0008727568907153145
$ws.Range("T10").Value = 0.0009066473515261946

$ws.Range("G11").Value = 163.8590903333333
$ws.Range("H11").Value = 491.577271
$ws.Range("I11").Value = 0.2750214756820535
$ws.Range("J11").Value = 0.284084617144743
$ws.Range("M11").Value = 28.36841766666667
$ws.Range("N11").Value = 85.105253
$ws.Range("O11").Value = 0.377369322353974
$ws.Range("P11").Value = 0.3795164395660363
$ws.Range("Q11").Value = 4648.423113056063
$ws.Range("R11").Value = 41835.80801750456
$ws.Range("S11").Value = 0.1037846679109265
$ws.Range("T11").Value = 0.1078147824342534

$ws.Range("G12").Value = 163.8590903333333
$ws.Range("H12").Value = 491.577271
$ws.Range("I12").Value = 0.2750214756820535
$ws.Range("J12").Value = 0.284084617144743
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.5
$ws.Range("M12").Value = 1.275895
$ws.Range("N12").Value = 2.55179
$ws.Range("O12").Value = 0.01697252336039082
$ws.Range("P12").Value = 0.01137939458707931
$ws.Range("Q12").Value = 209.0669940608483
$ws.Range("R12").Value = 1254.40196436509
$ws.Range("S12").Value = 0.004667808420622808
$ws.Range("T12").Value = 0.003232710954609387

$ws.Range("G13").Value = 163.8590903333333
$ws.Range("H13").Value = 491.577271
$ws.Range("I13").Value = 0.2750214756820535
$ws.Range("J13").Value = 0.284084617144743
$ws.Range("M13").Value = 45.29127766666667
$ws.Range("N13").Value = 135.873833
$ws.Range("O13").Value = 0.6024847406874758
$ws.Range("P13").Value = 0.6059126964859642
$ws.Range("Q13").Value = 7421.387558494416
$ws.Range("R13").Value = 66792.48802644973
$ws.Range("S13").Value = 0.1656962424597889
$ws.Range("T13").Value = 0.172130476404354

$ws.Range("G14").Value = 57.0238095
$ws.Range("H14").Value = 114.047619
$ws.Range("I14").Value = 0.09570889357312636
$ws.Range("J14").Value = 0.06590860906562239
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.2385586666666667
$ws.Range("N14").Value = 0.715676
$ws.Range("O14").Value = 0.003173413598159477
$ws.Range("P14").Value = 0.003191469360920207
$ws.Range("Q14").Value = 13.603523962574
$ws.Range("R14").Value = 81.621143775444
$ws.Range("S14").Value = 0.0003037239043297574
$ws.Range("T14").Value = 0.0002103453064538016

$ws.Range("G15").Value = 57.0238095
$ws.Range("H15").Value = 114.047619
$ws.Range("I15").Value = 0.09570889357312636
$ws.Range("J15").Value = 0.06590860906562239
$ws.Range("M15").Value = 28.36841766666667
$ws.Range("N15").Value = 85.105253
$ws.Range("O15").Value = 0.377369322353974
$ws.Range("P15").Value = 0.3795164395660363
$ws.Range("Q15").Value = 1617.675244840435
$ws.Range("R15").Value = 9706.051469042608
$ws.Range("S15").Value = 0.03611760031093931
$ws.Range("T15").Value = 0.02501340064933479

$ws.Range("G16").Value = 57.0238095
$ws.Range("H16").Value = 114.047619
$ws.Range("I16").Value = 0.09570889357312636
$ws.Range("J16").Value = 0.06590860906562239
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.5
$ws.Range("M16").Value = 1.275895
$ws.Range("N16").Value = 2.55179
$ws.Range("O16").Value = 0.01697252336039082
$ws.Range("P16").Value = 0.01137939458707931
$ws.Range("Q16").Value = 72.7563934220025
$ws.Range("R16").Value = 291.02557368801
$ws.Range("S16").Value = 0.001624421431967046
$ws.Range("T16").Value = 0.00075000006924327

$ws.Range("G17").Value = 57.0238095
$ws.Range("H17").Value = 114.047619
$ws.Range("I17").Value = 0.09570889357312636
$ws.Range("J17").Value = 0.06590860906562239
$ws.Range("M17").Value = 45.29127766666667
$ws.Range("N17").Value = 135.873833
$ws.Range("O17").Value = 0.6024847406874758
$ws.Range("P17").Value = 0.6059126964859642
$ws.Range("Q17").Value = 2582.681189675604
$ws.Range("R17").Value = 15496.08713805363
$ws.Range("S17").Value = 0.05766314792589025
$ws.Range("T17").Value = 0.03993486304059053

$ws.Range("G18").Value = 147.8896333333333
$ws.Range("H18").Value = 443.6689
$ws.Range("I18").Value = 0.2482183021684772
$ws.Range("J18").Value = 0.2563981636887546
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.2385586666666667
$ws.Range("N18").Value = 0.715676
$ws.Range("O18").Value = 0.003173413598159477
$ws.Range("P18").Value = 0.003191469360920207
$ws.Range("Q18").Value = 35.28035374182222
$ws.Range("R18").Value = 317.5231836764
$ws.Range("S18").Value = 0.0007876993354135036
$ws.Range("T18").Value = 0.0008182868836088641

$ws.Range("G19").Value = 147.8896333333333
$ws.Range("H19").Value = 443.6689
$ws.Range("I19").Value = 0.2482183021684772
$ws.Range("J19").Value = 0.2563981636887546
$ws.Range("M19").Value = 28.36841766666667
$ws.Range("N19").Value = 85.105253
$ws.Range("O19").Value = 0.377369322353974
$ws.Range("P19").Value = 0.3795164395660363
$ws.Range("Q19").Value = 4195.394886970189
$ws.Range("R19").Value = 37758.5539827317
$ws.Range("S19").Value = 0.09366997248517218
$ws.Range("T19").Value = 0.0973073181944259

$ws.Range("G20").Value = 147.8896333333333
$ws.Range("H20").Value = 443.6689
$ws.Range("I20").Value = 0.2482183021684772
$ws.Range("J20").Value = 0.2563981636887546
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.5
$ws.Range("M20").Value = 1.275895
$ws.Range("N20").Value = 2.55179
$ws.Range("O20").Value = 0.01697252336039082
$ws.Range("P20").Value = 0.01137939458707931
$ws.Range("Q20").Value = 188.6916437218333
$ws.Range("R20").Value = 1132.149862331
$ws.Range("S20").Value = 0.004212890932031027
$ws.Range("T20").Value = 0.00291765587601689

$ws.Range("G21").Value = 147.8896333333333
$ws.Range("H21").Value = 443.6689
$ws.Range("I21").Value = 0.2482183021684772
$ws.Range("J21").Value = 0.2563981636887546
$ws.Range("M21").Value = 45.29127766666667
$ws.Range("N21").Value = 135.873833
$ws.Range("O21").Value = 0.6024847406874758
$ws.Range("P21").Value = 0.6059126964859642
$ws.Range("Q21").Value = 6698.110447321522
$ws.Range("R21").Value = 60282.9940258937
$ws.Range("S21").Value = 0.1495477394158605
$ws.Range("T21").Value = 0.1553549027347029
